# Singapore Premier League workbook update (01-06-2024 01:16 base refresh)
#
# The underlying data refresh re-sorted the match rows for a handful of
# fixtures that share the same kickoff date, which manifests as the full
# row content (every column except the running index in column A) being
# swapped between two adjacent rows for each affected pair.
#
# NOTE: reading a Range/Cells ".Value" and feeding it straight back into
# another ".Value" assignment does not round-trip scalars correctly in
# this host - use ".Value2" for the read side (plain scalar) and ".Value"
# for the write side.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose full contents (columns B..AD) are swapped.
$pairs = @(
    @(4, 5),
    @(6, 7),
    @(22, 23),
    @(30, 31),
    @(54, 55)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
